# Add OOD detection with MSP and Energy score
# - Delete the extra "FPM" row (old row 11) from the DM aggregation group,
#   which shifts every row below it up by one.
# - Fill in the freshly-computed AUROC/FPR/ERR/AUPR results for the
#   DM score (Softmax / FL+1 aggregations).
# - Rename the leftover "FPM" aggregation label (now on the TRUSTED row)
#   to "PM".
# - Update the active selection to reflect where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "FPM" row under the DM score group (old row 11).
$ws.Rows(11).Delete()

# Fill in the new DM / FL results (row 9).
$ws.Range("E9").Value = 87.4
$ws.Range("F9").Value = 74.099999999999994
$ws.Range("G9").Value = 15.1
$ws.Range("H9").Value = 18.7

# Fill in the new DM / FL+1 results (row 10).
$ws.Range("E10").Value = 76.400000000000006
$ws.Range("F10").Value = 65.400000000000006
$ws.Range("G10").Value = 9.0500000000000007
$ws.Range("H10").Value = 20.9

# The former "FPM" aggregation (now on the TRUSTED row after the row
# shift above) becomes "PM".
$ws.Range("D13").Value = "PM"

# Match the author's final selection.
$ws.Range("G18").Select()
